$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data1 = @(
    @(1.02,1.027962381949805,1.045231825614765,1.027988166594529,1.050191063879345),
    @(1.02,1.02882584588133,1.04584376820937,1.028719369715356,1.050999447347006),
    @(1.02,1.029385229942504,1.046240226850495,1.02919346116083,1.051523515572686),
    @(1.02,1.02962055273454,1.046407013703913,1.029392995884252,1.051744068839513),
    @(1.02,1.029660073672571,1.046435024677215,1.029426511884272,1.051781114423807),
    @(1.02,1.029388373719268,1.046242455012564,1.029196126465063,1.051526461695465),
    @(1.02,1.028254054994801,1.045438531718014,1.028235081506536,1.050464054631241),
    @(1.02,1.026260401585669,1.0440257583497,1.02654898201188,1.048599630571652),
    @(1.02,1.024934859074966,1.04308660457672,1.025429984514258,1.047361961868533),
    @(1.02,1.02436174767518,1.042680603568425,1.024946670456541,1.046827315935871),
    @(1.02,1.024148998877867,1.042529897615441,1.024767331054434,1.046628917946388),
    @(1.02,1.024194628296833,1.042562219957179,1.024805791566624,1.046671466207094),
    @(1.02,1.02434415912656,1.042668144097912,1.024931842422225,1.046810912340497),
    @(1.02,1.024436307316157,1.042733420884546,1.02500953107661,1.046896855384311),
    @(1.02,1.024972912780505,1.043113563570533,1.025462086316614,1.047397471553095),
    @(1.02,1.025309741781299,1.043352194766747,1.025746289945936,1.047711837088263),
    @(1.02,1.025506291050458,1.043491447739254,1.02591217863447,1.047895323851049),
    @(1.02,1.025573323238716,1.043538940104982,1.025968762269524,1.047957908837794),
    @(1.02,1.025273594661478,1.043326585322705,1.025715785417701,1.047678095938571),
    @(1.02,1.02430012241532,1.042636949277812,1.024894718466929,1.046769843576564),
    @(1.02,1.023688816366813,1.042203932444195,1.024379552080431,1.046199908764751),
    @(1.02,1.024012809240415,1.042433426793932,1.024652549416982,1.046501935177856),
    @(1.02,1.025289927733905,1.043338156932379,1.02572956873807,1.047693341733005),
    @(1.02,1.026775187932826,1.044390527272205,1.026983993334426,1.049080707043333)
)

$data2 = @(
    @(1.037568581012903,1.033117341719321,1.04800064490111,1.030806292901606,1.052946027730943,1.034584486844778),
    @(1.037743714415541,1.033621654474429,1.048424505975594,1.031345543910449,1.053566818494378,1.035089515781843),
    @(1.037856141267098,1.033948004333404,1.04869856794312,1.031694761646936,1.053968846536131,1.035416329095051),
    @(1.03790319025516,1.034085206759929,1.048813733275878,1.031841639991947,1.054137937605476,1.035553726364759),
    @(1.037911077343549,1.034108243905434,1.048833067051665,1.031866305414958,1.054166333292853,1.035576796225653),
    @(1.037856770784697,1.033949837621047,1.048700106987362,1.031696723980152,1.053971105633033,1.035418164986172),
    @(1.037627953293901,1.033287770499425,1.048143932257547,1.030988475232242,1.053155756511298,1.034755157653299),
    @(1.037217915729949,1.03212137539312,1.04716237544088,1.029742701625304,1.051721647212633,1.033587106132155),
    @(1.03694000673516,1.031344021042801,1.046507073165695,1.028913778528415,1.050767450685119,1.032808647849495),
    @(1.036818598747012,1.031007490918318,1.046223116423562,1.028555240107645,1.050354739521429,1.032471639813635),
    @(1.036773342138919,1.030882500313437,1.046117612605252,1.028422123024123,1.050201511397276,1.032346471707757),
    @(1.036783057076308,1.030909310674665,1.046140244851464,1.028450674351452,1.050234376101276,1.032373320142774),
    @(1.036814861087107,1.030997158910772,1.046214396041674,1.028544235375222,1.050342072174854,1.032461293133452),
    @(1.036834435371141,1.031051286691388,1.04626007915103,1.028601889407265,1.050408436763813,1.032515497781725),
    @(1.036948041667874,1.031366357007219,1.046525914206581,1.028937581889972,1.050794850852432,1.032831015533544),
    @(1.037019017513039,1.03156401175458,1.046692611199153,1.02914825847926,1.051037363177038,1.033028950973317),
    @(1.037060313084774,1.031679307041647,1.046789822714178,1.029271180181021,1.051178860995236,1.033144409992917),
    @(1.037074376240993,1.031718620839838,1.046822965890574,1.029313099643428,1.051227115586828,1.033183779621211),
    @(1.037011413170485,1.031542804577436,1.046674728250252,1.029125650975009,1.051011339286212,1.033007713679549),
    @(1.036805500017269,1.030971289462801,1.046192561171423,1.028516682317205,1.050310356386662,1.032435386947898),
    @(1.036675107294662,1.030612023342537,1.045889232783986,1.028134147647524,1.049870033287208,1.032075610628535),
    @(1.036744318569333,1.030802470241249,1.046050048628655,1.028336903069496,1.050103417198371,1.032266327983688),
    @(1.037014849565254,1.031552387167639,1.046682808847136,1.029135866217005,1.051023098225071,1.03301730987813),
    @(1.03732472548448,1.032422879121261,1.047416301470133,1.030064489039854,1.052092075541583,1.033889038030173)
)

$startRow = 2
$cols1 = @(2,3,4,5,6)
$cols2 = @(9,10,11,12,13,14)

for ($i = 0; $i -lt $data1.Count; $i++) {
    $r = $startRow + $i
    $rowVals1 = $data1[$i]
    for ($j = 0; $j -lt $cols1.Count; $j++) {
        $ws.Cells.Item($r, $cols1[$j]).Value = $rowVals1[$j]
    }
    $rowVals2 = $data2[$i]
    for ($j = 0; $j -lt $cols2.Count; $j++) {
        $ws.Cells.Item($r, $cols2[$j]).Value = $rowVals2[$j]
    }
}
